$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 337, pushing the existing 337:391 block down to 338:392.
$ws.Rows.Item(337).Insert()

# Populate the newly inserted row with the new weekly record. All the
# "static" columns (A,B,C,E,F,G,H,I,N,O,Q,R) repeat the same values used
# throughout this market/variety block; only the per-record columns
# (D,J,K,L,M,P) carry the new observation's data.
$ws.Cells.Item(337, 1).Value = 4
$ws.Cells.Item(337, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(337, 3).Value = "Los Lagos"
$ws.Cells.Item(337, 4).Value = 44951
$ws.Cells.Item(337, 5).Value = 10
$ws.Cells.Item(337, 6).Value = 100112037
$ws.Cells.Item(337, 7).Value = "Cebollín"
$ws.Cells.Item(337, 8).Value = "Sin especificar"
$ws.Cells.Item(337, 9).Value = "Primera"
$ws.Cells.Item(337, 10).Value = 35
$ws.Cells.Item(337, 11).Value = 6000
$ws.Cells.Item(337, 12).Value = 6000
$ws.Cells.Item(337, 13).Value = 6000
$ws.Cells.Item(337, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(337, 15).Value = "Región Metropolitana"
$ws.Cells.Item(337, 16).Value = 167
$ws.Cells.Item(337, 17).Value = 36
$ws.Cells.Item(337, 18).Value = "Hortaliza"
